$d = $word.ActiveDocument

# 1) Table 1, Row 1, Col 2: "4주차" -> "3주차" (change only the leading digit run)
$cell1 = $d.Tables.Item(1).Cell(1, 2)
$r1 = $cell1.Range
$r1.End = $r1.Start + 1
$r1.Text = "3"

# 2) Table 2, Row 3, Col 2: "5주차" -> "4주차" (change only the leading digit run)
$cell2 = $d.Tables.Item(2).Cell(3, 2)
$r2 = $cell2.Range
$r2.End = $r2.Start + 1
$r2.Text = "4"
